$wb = $excel.ActiveWorkbook

# --- Sheet "times": update StartTime/StopTime (B2/B3) ---
$wsTimes = $wb.Worksheets.Item("times")
$wsTimes.Range("B2").Value = 43830.99861111111
$wsTimes.Range("B3").Value = 44195.99861111111

# --- Sheet "scenario_data_emlab": drop the 2025 column (C), update 2020 (B) values ---
$wsScenario = $wb.Worksheets.Item("scenario_data_emlab")
$wsScenario.Range("C1:C8").Delete()

$wsScenario.Range("B1").Value = 2020
$wsScenario.Range("B2").Value = 20.4
$wsScenario.Range("B5").Value = 10.8
$wsScenario.Range("B6").Value = 20.16
$wsScenario.Range("B7").Value = 46.44

# --- Sheet "renewables": drop the bulk placeholder rows, keep only 4 plants ---
$wsRenew = $wb.Worksheets.Item("renewables")
$wsRenew.Rows("6:27").Delete()

$wsRenew.Range("B2").Value = 20152400023
$wsRenew.Range("C2").Value = 47547.50848700004
$wsRenew.Range("D2").Value = 1.35
$wsRenew.Range("E2").Value = "WindOn"

$wsRenew.Range("B3").Value = 20151200026
$wsRenew.Range("C3").Value = 8858.749999999998
$wsRenew.Range("D3").Value = 0
$wsRenew.Range("E3").Value = "RunOfRiver"

$wsRenew.Range("B4").Value = 20152100030
$wsRenew.Range("C4").Value = 53555.51607579708
$wsRenew.Range("D4").Value = 0
$wsRenew.Range("E4").Value = "OtherPV"

$wsRenew.Range("B5").Value = 20152300031
$wsRenew.Range("C5").Value = 10271.8
$wsRenew.Range("D5").Value = 2.7
$wsRenew.Range("E5").Value = "WindOff"
